{"js": "// Replace every occurrence of the old \"Datas das campanhas...\" sentence\n// with the reworded version throughout the document body.\nconst oldText =\n  \"Datas das campanhas de Constela\u00e7\u00e3o de Cygnus 2022: 10 a 19 de agosto, 9 a 18 de setembro, 8 a 17 de outubro\";\nconst newText =\n  \"Datas das campanhas de 2022 que usam Constela\u00e7\u00e3o de Cygnus: 10 a 19 de agosto, 9 a 18 de setembro, 8 a 17 de outubro\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace every occurrence of the old \"Datas das campanhas...\" sentence\n# with the reworded version throughout the document.\n$d = $word.ActiveDocument\n\n$oldText = \"Datas das campanhas de Constela\u00e7\u00e3o de Cygnus 2022: 10 a 19 de agosto, 9 a 18 de setembro, 8 a 17 de outubro\"\n$newText = \"Datas das campanhas de 2022 que usam Constela\u00e7\u00e3o de Cygnus: 10 a 19 de agosto, 9 a 18 de setembro, 8 a 17 de outubro\"\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n$find.Execute(\n    $find.Text,      # FindText\n    $false,          # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    $wdFindContinue, # Wrap\n    $false,          # Format\n    $find.Replacement.Text, # ReplaceWith\n    $wdReplaceAll    # Replace\n) | Out-Null\n"}
